# Update leve market/profit data across sheets (scheduled data refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 18: You Grow, Girl / Growth Formula Beta
$ws.Range("H18").Value = 3500
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

# Row 33: Glazed and Confused / Clear Glass Lens
$ws.Range("H33").Value = 326.70587
$ws.Range("I33").Value = 111.888885
$ws.Range("K33").Value = 111.888885
$ws.Range("M33").Value = 117.111115

# Row 43: Growing Is Knowing / Growth Formula Gamma
$ws.Range("H43").Value = 2393.4546
$ws.Range("I43").Value = 2148
$ws.Range("J43").Value = 2598
$ws.Range("K43").Value = 2148
$ws.Range("L43").Value = 2598
$ws.Range("M43").Value = -2079
$ws.Range("N43").Value = -2736

# Row 94: Magic Beans / Growth Formula Eta
$ws.Range("H94").Value = 1800
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()

# Row 135: For Tired Minds / Grade 1 Gemsap of Intelligence
$ws.Range("H135").Value = 582.6667
$ws.Range("I135").Value = 617.4545000000001
$ws.Range("K135").Value = 5557.0905
$ws.Range("M135").Value = -3022.0905

$ws = $wb.Worksheets.Item("ARM")
# Row 5: The Alloyed Truth / Bronze Rivets
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("N5").ClearContents()

# Row 13: Get into Their Heads / Bronze Chain Coif
$ws.Range("H13").Value = 17498.2
$ws.Range("I13").Value = 15000
$ws.Range("K13").Value = 15000
$ws.Range("M13").Value = -14856

# Row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 12976.238
$ws.Range("I32").Value = 8115.1333
$ws.Range("K32").Value = 8115.1333
$ws.Range("M32").Value = -7828.1333

# Row 45: Hollow Hallmarks / Mythril Ingot
$ws.Range("H45").Value = 1656.9
$ws.Range("I45").Value = 1696.125
$ws.Range("K45").Value = 1696.125
$ws.Range("M45").Value = -1319.125

# Row 102: Smells of Rich Tama-hagane / Tama-hagane Ingot
$ws.Range("H102").Value = 1685.0769
$ws.Range("I102").Value = 1685.0769
$ws.Range("K102").Value = 1685.0769
$ws.Range("M102").Value = -63.07690000000002

# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 774.2727
$ws.Range("I132").Value = 774.2727
$ws.Range("K132").Value = 2322.8181
$ws.Range("M132").Value = 207.1819

$ws = $wb.Worksheets.Item("BSM")
# Row 4: Mending Fences / Bronze Rivets
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").ClearContents()

# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 3352.8845
$ws.Range("I134").Value = 3352.8845
$ws.Range("K134").Value = 10058.6535
$ws.Range("M134").Value = -7523.6535

$ws = $wb.Worksheets.Item("CRP")
# Row 22: Driving Up the Wall / Elm Lumber
$ws.Range("H22").Value = 679.875
$ws.Range("I22").Value = 656.6667
$ws.Range("J22").Value = 749.5
$ws.Range("K22").Value = 656.6667
$ws.Range("L22").Value = 749.5
$ws.Range("M22").Value = -306.6667
$ws.Range("N22").Value = -1449.5

# Row 107: Built to Last / White Oak Lumber
$ws.Range("H107").Value = 936.55554
$ws.Range("I107").Value = 996.7143
$ws.Range("K107").Value = 996.7143
$ws.Range("M107").Value = 923.2857

# Row 122: Timber of Tenkonto / Horse Chestnut Lumber
$ws.Range("H122").Value = 1274.875
$ws.Range("I122").Value = 1033.1666
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 3099.4998
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -649.4998000000001
$ws.Range("N122").Value = -10900

$ws = $wb.Worksheets.Item("CUL")
# Row 5: What a Sap / Maple Syrup
$ws.Range("H5").Value = 1371.6154
$ws.Range("J5").Value = 472.66666
$ws.Range("L5").Value = 1417.99998
$ws.Range("N5").Value = -1641.99998

# Row 14: Keep Your Powder Dry / Kukuru Powder
$ws.Range("H14").Value = 16981.834
$ws.Range("I14").Value = 16981.834
$ws.Range("K14").Value = 50945.50199999999
$ws.Range("M14").Value = -50772.50199999999

# Row 50: Moving Up in the World / Rolanberry Cheese
$ws.Range("H50").Value = 555
$ws.Range("I50").Value = 555
$ws.Range("K50").Value = 1665
$ws.Range("M50").Value = -1184

# Row 53: Rolanberry Fields Forever / Rolanberry Cheese
$ws.Range("H53").Value = 555
$ws.Range("I53").Value = 555
$ws.Range("K53").Value = 1665
$ws.Range("M53").Value = -1184

# Row 80: Saucy for a Suitor / Hollandaise Sauce
$ws.Range("H80").Value = 6284.4287
$ws.Range("J80").Value = 10623.75
$ws.Range("L80").Value = 31871.25
$ws.Range("N80").Value = -33743.25

# Row 83: Saved by the Sauce (L) / Hollandaise Sauce
$ws.Range("H83").Value = 6284.4287
$ws.Range("J83").Value = 10623.75
$ws.Range("L83").Value = 95613.75
$ws.Range("N83").Value = -104973.75

# Row 98: Sweet Kiss of Death / Rice Vinegar
$ws.Range("H98").Value = 373.33334
$ws.Range("J98").Value = 372.5
$ws.Range("L98").Value = 1117.5
$ws.Range("N98").Value = -4113.5

# Row 107: Slippery Service / Frantoio Oil
$ws.Range("H107").Value = 512
$ws.Range("I107").Value = 383.45456
$ws.Range("K107").Value = 1150.36368
$ws.Range("M107").Value = 769.6363200000001

# Row 135: Not-so-secret Ingredient / Royal Maple Syrup
$ws.Range("H135").Value = 1371.6154
$ws.Range("J135").Value = 472.66666
$ws.Range("L135").Value = 4253.99994
$ws.Range("N135").Value = -9323.99994

$ws = $wb.Worksheets.Item("GSM")
# Row 97: If I'd a Koppranickel for Every Time... / Koppranickel Ingot
$ws.Range("H97").Value = 656.5294
$ws.Range("I97").Value = 539.25
$ws.Range("J97").Value = 938
$ws.Range("K97").Value = 539.25
$ws.Range("L97").Value = 938
$ws.Range("M97").Value = -43.25
$ws.Range("N97").Value = -1930

$ws = $wb.Worksheets.Item("LTW")
# Row 22: Skin off Their Backs / Aldgoat Leather
$ws.Range("H22").Value = 1511.4849
$ws.Range("I22").Value = 1513.32
$ws.Range("J22").Value = 1505.75
$ws.Range("K22").Value = 1513.32
$ws.Range("L22").Value = 1505.75
$ws.Range("M22").Value = -1218.32
$ws.Range("N22").Value = -2095.75

# Row 27: Fire and Hide / Aldgoat Leather
$ws.Range("H27").Value = 1511.4849
$ws.Range("I27").Value = 1513.32
$ws.Range("J27").Value = 1505.75
$ws.Range("K27").Value = 1513.32
$ws.Range("L27").Value = 1505.75
$ws.Range("M27").Value = -1406.32
$ws.Range("N27").Value = -1719.75

# Row 46: Supply Side Logic / Boar Leather
$ws.Range("H46").Value = 2400
$ws.Range("I46").Value = 2600
$ws.Range("J46").Value = 2000
$ws.Range("K46").Value = 2600
$ws.Range("L46").Value = 2000
$ws.Range("M46").Value = -2412
$ws.Range("N46").Value = -2376

# Row 55: It's Not a Job, It's a Calling / Peiste Leather
$ws.Range("H55").Value = 5162.25
$ws.Range("J55").Value = 225
$ws.Range("L55").Value = 225
$ws.Range("N55").Value = -571

# Row 74: Overall, We Blend In / Dhalmelskin Vest
$ws.Range("H74").Value = 38331.332
$ws.Range("I74").Value = 15000
$ws.Range("K74").Value = 15000
$ws.Range("M74").Value = -14002

# Row 77: Eviction Notice (L) / Dhalmelskin Vest
$ws.Range("H77").Value = 38331.332
$ws.Range("I77").Value = 15000
$ws.Range("K77").Value = 45000
$ws.Range("M77").Value = -40008

$ws = $wb.Worksheets.Item("WVR")
# Row 10: Just for Kecks / Hempen Kecks
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()

# Row 136: Weaving the Envelope / Sarcenet Cloth
$ws.Range("H136").Value = 1075.7858
$ws.Range("I136").Value = 1075.7858
$ws.Range("K136").Value = 3227.3574
$ws.Range("M136").Value = -677.3574000000003
